$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: force a cell value to be stored as text (avoids Excel auto-number conversion)
# while leaving the cells style/number-format unchanged afterwards.
function Set-TextValue($addr, $val) {
    $r = $ws.Range($addr)
    $r.NumberFormat = "@"
    $r.Value = $val
    $r.Style = "Normal"
}

$ws.Range("D2").Value = "27.154.84"
$ws.Range("E2").Value = "  +3.53%  "
$ws.Range("D3").Value = "1.661.12"
$ws.Range("E3").Value = "  +4.03%  "
$ws.Range("E4").Value = "  -0.02%  "
Set-TextValue "D5" "215.64"
$ws.Range("E5").Value = "  +1.81%  "
$ws.Range("E6").Value = "  +1.11%  "
$ws.Range("E7").Value = "  -0.09%  "
$ws.Range("E8").Value = "  +2.37%  "
Set-TextValue "D10" "19.65"
$ws.Range("E10").Value = "  +3.81%  "
$ws.Range("E11").Value = "  +1.22%  "
$ws.Range("D12").Value = "1.894.89"
$ws.Range("E12").Value = "  +4.01%  "
$ws.Range("D13").Value = "1.661.98"
$ws.Range("E13").Value = "  +4.14%  "
$ws.Range("E14").Value = "  +1.90%  "
$ws.Range("E15").Value = "  +3.26%  "
Set-TextValue "D16" "64.96"
Set-TextValue "D17" "240.77"
$ws.Range("E17").Value = "  +5.68%  "
$ws.Range("D18").Value = "27.131.27"
$ws.Range("E18").Value = "  +3.48%  "
Set-TextValue "D19" "7.87"
$ws.Range("E19").Value = "  +4.02%  "
$ws.Range("E20").Value = "  +1.60%  "
$ws.Range("E21").Value = "  -0.05%  "
Set-TextValue "D22" "4.45"
$ws.Range("E22").Value = "  +5.16%  "
Set-TextValue "D23" "2.25"
$ws.Range("E23").Value = "  +4.04%  "
Set-TextValue "D24" "9.34"
$ws.Range("E24").Value = "  +5.11%  "
Set-TextValue "D25" "145.80"
$ws.Range("E25").Value = "  +0.22%  "
$ws.Range("E26").Value = "  -0.05%  "
$ws.Range("E27").Value = "  +3.15%  "
$ws.Range("E28").Value = "  +1.25%  "
$ws.Range("E29").Value = "  +3.40%  "
$ws.Range("E30").Value = "  +1.18%  "
$ws.Range("E31").Value = "  +1.39%  "
$ws.Range("D32").Value = "1.539.23"
$ws.Range("E32").Value = "  +6.29%  "
$ws.Range("E33").Value = "  +3.21%  "
$ws.Range("E34").Value = "  +4.01%  "
$ws.Range("E35").Value = "  +8.57%  "
$ws.Range("E36").Value = "  +0.05%  "
$ws.Range("E37").Value = "  +2.21%  "
Set-TextValue "D38" "0.896"
$ws.Range("E38").Value = "  +9.65%  "
$ws.Range("E39").Value = "  +3.31%  "
Set-TextValue "D40" "5.97"
$ws.Range("E40").Value = "  +4.20%  "
$ws.Range("E41").Value = "  -0.09%  "
$ws.Range("E42").Value = "  +4.55%  "
Set-TextValue "D43" "66.22"
$ws.Range("E43").Value = "  +9.70%  "
$ws.Range("D44").Value = "1.801.06"
$ws.Range("E44").Value = "  +3.81%  "
$ws.Range("E45").Value = "  +2.23%  "
$ws.Range("E46").Value = "  -1.12%  "
Set-TextValue "D47" "90.48"
$ws.Range("E47").Value = "  +3.33%  "

# Row 48/49 swap: BabyDogeCoin moves to row 48 (was RenderToken), RenderToken moves to row 49 (was BabyDogeCoin)
$ws.Range("B48").Value = "BabyDogeCoin"
$ws.Range("C48").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
Set-TextValue "D48" "0.0₆0105"
$ws.Range("E48").Value = "  +0.34%  "

$ws.Range("B49").Value = "RenderToken"
$ws.Range("C49").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
Set-TextValue "D49" "1.54"
$ws.Range("E49").Value = "  +4.20%  "
$ws.Range("E50").Value = "  +0.82%  "
$ws.Range("E51").Value = "  +3.60%  "
